# Weekly refresh: a new "Poroto granado" price-report row for
# Terminal La Palmera de La Serena is inserted at row 44, pushing the
# previously existing rows 44-68 down to 45-69 (dimension grows to R69).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 44 (shifts 44..68 -> 45..69,
# and carries the existing formatting, e.g. the date style on column D).
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with this week's data.
$ws.Cells.Item(44, 1).Value  = 8
$ws.Cells.Item(44, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(44, 3).Value  = "Coquimbo"
$ws.Cells.Item(44, 4).Value  = 44603
$ws.Cells.Item(44, 5).Value  = 4
$ws.Cells.Item(44, 6).Value  = 100112030
$ws.Cells.Item(44, 7).Value  = "Poroto granado"
$ws.Cells.Item(44, 8).Value  = "Sin especificar"
$ws.Cells.Item(44, 9).Value  = "Primera"
$ws.Cells.Item(44, 10).Value = 520
$ws.Cells.Item(44, 11).Value = 31000
$ws.Cells.Item(44, 12).Value = 32000
$ws.Cells.Item(44, 13).Value = 31500
$ws.Cells.Item(44, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(44, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(44, 16).Value = 1260
$ws.Cells.Item(44, 17).Value = 25
$ws.Cells.Item(44, 18).Value = "Hortaliza"
